# Se obtiene id del cliente a través del DNI.
# Insert a new "ClienteDNI" column right before the "Cliente" column on both
# sheets ("Pendientes" and "Facturados"), populating sample DNI values on
# the "Pendientes" data rows.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Pendientes" --------------------------------------------------
$ws1 = $wb.Worksheets.Item("Pendientes")

# Column F ("Cliente" was F, now becomes G) - insert a fresh column at F.
$ws1.Columns.Item(6).Insert()

# Header
$ws1.Range("F1").Value() = "ClienteDNI"

# Give the new data cells the same look (number format / alignment) as the
# neighbouring "FORMA DE PAGO"/"FACT/REMIS"/"DCTO FNCIERO" cells in column B.
$ws1.Range("B2:B3").Copy()
$ws1.Range("F2:F3").PasteSpecial(-4122)
$ws1.Range("F2:F3").Value() = 12345

# Restore the shared formula in column S (was R before the insert).
$ws1.Range("S2:S3").Formula() = "=Q2*R2"

[void]$ws1.Range("F1").Select()

# ---- Sheet "Facturados" --------------------------------------------------
$ws2 = $wb.Worksheets.Item("Facturados")
$ws2.Columns.Item(6).Insert()
$ws2.Range("F1").Value() = "ClienteDNI"

[void]$ws2.Range("F5").Select()
